$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.072.23"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "2.629.00"
$ws.Range("E3").Value = "  -2.80%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'594.06"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").Value = "'166.57"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.533"
$ws.Range("E8").Value = "  -3.30%  "
$ws.Range("D9").Value = "2.631.02"
$ws.Range("E9").Value = "  -2.84%  "
$ws.Range("D10").Value = "'0.141"
$ws.Range("E10").Value = "  -2.28%  "
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D12").Value = "'0.360"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").Value = "'27.64"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").Value = "3.115.16"
$ws.Range("E15").Value = "  -2.74%  "
$ws.Range("D16").Value = "'0.0000181"
$ws.Range("E16").Value = "  -2.96%  "
$ws.Range("D17").Value = "67.056.63"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "2.628.15"
$ws.Range("E18").Value = "  -3.00%  "
$ws.Range("D19").Value = "'12.22"
$ws.Range("E19").Value = "  +4.15%  "
$ws.Range("D20").Value = "'8.09"
$ws.Range("E20").Value = "  +6.29%  "
$ws.Range("D21").Value = "'357.64"
$ws.Range("E21").Value = "  -3.11%  "
$ws.Range("D22").Value = "'4.33"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").Value = "'4.66"
$ws.Range("E23").Value = "  -5.20%  "
$ws.Range("D24").Value = "'10.93"
$ws.Range("E24").Value = "  +10.20%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("B26").Value = "SuiNetwork"
$ws.Range("C26").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D26").Value = "'1.94"
$ws.Range("E26").Value = "  -6.46%  "
$ws.Range("D27").Value = "'70.24"
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("D28").Value = "2.766.15"
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0000101"
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'552.82"
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'7.89"
$ws.Range("E32").Value = "  -2.67%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.36"
$ws.Range("E33").Value = "  -3.06%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "'1.90"
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.136"
$ws.Range("E35").Value = "  +3.81%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.50"
$ws.Range("E37").Value = "  -6.05%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'156.79"
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "'19.11"
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.365"
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("D41").Value = "'1.79"
$ws.Range("E41").Value = "  -3.40%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'5.15"
$ws.Range("E42").Value = "  -4.31%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'17.94"
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.45"
$ws.Range("E45").Value = "  -4.75%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'40.18"
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0298"
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.582"
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'151.63"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "'3.80"
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").Value = "'1.72"
$ws.Range("E51").Value = "  -2.11%  "
